$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate column A rows 6-36 with the new shared-string text values.
# Rows 1-5 already contain the original header/label text and are left untouched.
$ws.Range("A6").Value = "Data Center Environment"
$ws.Range("A7").Value = "Data Center Design Capacity"
$ws.Range("A8").Value = "Cooling System"
$ws.Range("A9").Value = "Air distribution type"
$ws.Range("A10").Value = "UPS Architecture"
$ws.Range("A11").Value = "Power distribution type"
$ws.Range("A12").Value = "Power Density"
$ws.Range("A13").Value = "Core & shell"
$ws.Range("A14").Value = "4 kW / rack"
$ws.Range("A15").Value = "$ 90 / hour"
$ws.Range("A16").Value = "$ 90 / ft²"
$ws.Range("A17").Value = "Labor Rate"
$ws.Range("A18").Value = "Redundancy Level"
$ws.Range("A19").Value = "Power"
$ws.Range("A20").Value = "Cooling"
$ws.Range("A21").Value = "IT distribution"
$ws.Range("A22").Value = "UPS"
$ws.Range("A23").Value = "Generator"
$ws.Range("A24").Value = "Capital Cost Summary"
$ws.Range("A25").Value = "$ 7.2 M"
$ws.Range("A26").Value = "Data Center Cost"
$ws.Range("A27").Value = "Data Center Cost Per Watt"
$ws.Range("A28").Value = "'$ 7.19"
$ws.Range("A29").Value = "Calculated Rack Quantity"
$ws.Range("A30").Value = "'250"
$ws.Range("A31").Value = "IT Room Area"
$ws.Range("A32").Value = "8,125 ft²"
$ws.Range("A33").Value = "Facility Area"
$ws.Range("A34").Value = "13,195 ft²"
$ws.Range("A35").Value = "Cost by Type"
$ws.Range("A36").Value = "Cost by System"

# Move the active cell/selection to reflect the final data-entry position.
$ws.Range("A22").Select()

